# Auto commit - 10171723
# Updates the HL Maintain Report: refresh the report-date stamp in the
# title cell, append one new maintenance record as row 74, mark the
# previously-last row's "work content" / P column with wrap formatting,
# and extend Print_Area / sheet dimension to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------------
# 1. Refresh the "製表日期" (report-generated-on) date in the title cell.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202510   (  製表日期:2025-10-17  )"

# ---------------------------------------------------------------------
# 2. Row 73 (previously the last data row) gains wrap-text formatting on
#    its P and AC columns (P73 stays blank; AC73 keeps its text).
# ---------------------------------------------------------------------
$ws.Range("P73").WrapText = $true
$ws.Range("AC73").WrapText = $true

# ---------------------------------------------------------------------
# 3. Build row 74 by cloning row 72's cell formatting (plain / non-
#    highlighted data row), then writing in the new record's values.
# ---------------------------------------------------------------------

# D74 ("客戶工作案號") is a long digit string that must stay TEXT, not get
# auto-coerced into a number. Stamp it as Text *before* the row-format
# clone below, so the subsequent format copy restores the normal
# (non-quote-prefixed) "General" look of the rest of the row.
$ws.Cells.Item(74, 4).NumberFormat = "@"
$ws.Cells.Item(74, 4).Value = "13770114101701"

$ws.Range("A72:AK72").Copy()
$ws.Range("A74:AK74").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# P/AC on this new row stay left-aligned without wrap (matches the rest
# of the non-highlighted rows, unlike row 73's special wrapped cells).
$ws.Range("P74").WrapText = $false
$ws.Range("AC74").WrapText = $false

$ws.Cells.Item(74, 1).Value = 72
$ws.Cells.Item(74, 2).Value = "維修"
$ws.Cells.Item(74, 3).Value = 2025102286
$ws.Cells.Item(74, 5).Value = "一般件"
$ws.Cells.Item(74, 6).Value = 3770
$ws.Cells.Item(74, 7).Value = "北縣西雲店"
$ws.Cells.Item(74, 8).Value = "新北市五股區"
$ws.Cells.Item(74, 9).Value = "2025-10-17 15:00:46"
$ws.Cells.Item(74, 10).Value = "星期五"
$ws.Cells.Item(74, 11).Value = "下午"
$ws.Cells.Item(74, 12).Value = "HL23"
$ws.Cells.Item(74, 13).Value = "HL-TM主機"
$ws.Cells.Item(74, 14).Value = 2307
$ws.Cells.Item(74, 15).Value = "觸控不良(游標偏移)"
$ws.Cells.Item(74, 16).Value = "門市告知TM2(TCX800)游標一直觸控下方，無法協助觸控校正重啟後仍異常，與門市確認無張貼文宣，10/13台芝到店清潔螢幕.觸控校正後有恢復正常(13770114100901)但隔天即又開始異常，門市店長表示如關機休息後重開就會正常，但如果按交班開機後就會又開始發生該情況，已造成門市作業不便.....須請台芝到店協助"
$ws.Cells.Item(74, 17).Value = "THILF03770"
$ws.Cells.Item(74, 18).Value = "新北一"
$ws.Cells.Item(74, 19).Value = "湯家瑋"
$ws.Cells.Item(74, 20).Value = 1
$ws.Cells.Item(74, 21).Value = "已完工"
$ws.Cells.Item(74, 22).Value = "2025-10-17 15:04:02"
$ws.Cells.Item(74, 23).Value = "2025-10-17 15:30:00"
$ws.Cells.Item(74, 24).Value = "2025-10-17 16:30:00"
$ws.Cells.Item(74, 25).Value = "2025-10-20 19:04:00"
$ws.Cells.Item(74, 26).Value = 1
$ws.Cells.Item(74, 28).Value = "到場處理"
$ws.Cells.Item(74, 29).Value = "更換客顯器 觸控校正"
$ws.Cells.Item(74, 37).Value = "O"

# ---------------------------------------------------------------------
# 4. Extend the print area to include the newly-added row.
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = "'Report'!`$A`$1:`$AK`$74"

# ---------------------------------------------------------------------
# 5. Move the selection cursor (matches where the editor last left off).
# ---------------------------------------------------------------------
$ws.Range("AC71").Select() | Out-Null
